# Swap the order of the first two comma-separated names/emails in the
# "Recorded By" column (column G) for every data row, except rows whose
# value is already exactly "System, dnasr281@gmail.com" (those are left
# untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    if ($text -eq "System, dnasr281@gmail.com") { continue }

    $parts = $text -split ",\s*"

    if ($parts.Count -ge 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()

        $swapped = @($second, $first)
        if ($parts.Count -gt 2) {
            $swapped += $parts[2..($parts.Count - 1)]
        }
        $newText = [string]::Join(", ", $swapped)

        $cell.Value2 = $newText
    }
}
